$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 38
$ws.Cells.Item(38, 4).Value = 44687
$ws.Cells.Item(38, 11).Value = "Sin especificar"
$ws.Cells.Item(38, 12).Value = "2a amarillo"
$ws.Cells.Item(38, 13).Value = 300
$ws.Cells.Item(38, 14).Value = 20000
$ws.Cells.Item(38, 15).Value = 22000
$ws.Cells.Item(38, 16).Value = 21000
$ws.Cells.Item(38, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(38, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(38, 19).Value = 1050
$ws.Cells.Item(38, 20).Value = 20

# Row 39
$ws.Cells.Item(39, 4).Value = 44363
$ws.Cells.Item(39, 11).Value = "Sin especificar"
$ws.Cells.Item(39, 12).Value = "2a amarillo"
$ws.Cells.Item(39, 13).Value = 200
$ws.Cells.Item(39, 14).Value = 11000
$ws.Cells.Item(39, 15).Value = 12000
$ws.Cells.Item(39, 16).Value = 11500
$ws.Cells.Item(39, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(39, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(39, 19).Value = 575
$ws.Cells.Item(39, 20).Value = 20

# Row 40
$ws.Cells.Item(40, 4).Value = 44363
$ws.Cells.Item(40, 11).Value = "Sutil De Gase"
$ws.Cells.Item(40, 12).Value = "Primera"
$ws.Cells.Item(40, 13).Value = 250
$ws.Cells.Item(40, 14).Value = 24000
$ws.Cells.Item(40, 15).Value = 25000
$ws.Cells.Item(40, 16).Value = 24500
$ws.Cells.Item(40, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(40, 18).Value = "Perú"
$ws.Cells.Item(40, 19).Value = 1021
$ws.Cells.Item(40, 20).Value = 24

# Row 41
$ws.Cells.Item(41, 4).Value = 44363
$ws.Cells.Item(41, 11).Value = "Tahití"
$ws.Cells.Item(41, 12).Value = "Primera"
$ws.Cells.Item(41, 13).Value = 270
$ws.Cells.Item(41, 14).Value = 25000
$ws.Cells.Item(41, 15).Value = 26000
$ws.Cells.Item(41, 16).Value = 25500
$ws.Cells.Item(41, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(41, 18).Value = "Perú"
$ws.Cells.Item(41, 19).Value = 1062
$ws.Cells.Item(41, 20).Value = 24

# Row 42
$ws.Cells.Item(42, 4).Value = 44244
$ws.Cells.Item(42, 11).Value = "Sin especificar"
$ws.Cells.Item(42, 12).Value = "1a amarillo"
$ws.Cells.Item(42, 13).Value = 260
$ws.Cells.Item(42, 14).Value = 29000
$ws.Cells.Item(42, 15).Value = 30000
$ws.Cells.Item(42, 16).Value = 29500
$ws.Cells.Item(42, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(42, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(42, 19).Value = 1475
$ws.Cells.Item(42, 20).Value = 20

# Row 43
$ws.Cells.Item(43, 4).Value = 44447
$ws.Cells.Item(43, 11).Value = "Sin especificar"
$ws.Cells.Item(43, 12).Value = "2a amarillo"
$ws.Cells.Item(43, 13).Value = 270
$ws.Cells.Item(43, 14).Value = 10000
$ws.Cells.Item(43, 15).Value = 11000
$ws.Cells.Item(43, 16).Value = 10500
$ws.Cells.Item(43, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(43, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(43, 19).Value = 525
$ws.Cells.Item(43, 20).Value = 20

# Row 44
$ws.Cells.Item(44, 4).Value = 44608
$ws.Cells.Item(44, 11).Value = "Tahití"
$ws.Cells.Item(44, 12).Value = "Primera"
$ws.Cells.Item(44, 13).Value = 300
$ws.Cells.Item(44, 14).Value = 34000
$ws.Cells.Item(44, 15).Value = 35000
$ws.Cells.Item(44, 16).Value = 34500
$ws.Cells.Item(44, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(44, 18).Value = "Perú"
$ws.Cells.Item(44, 19).Value = 1438
$ws.Cells.Item(44, 20).Value = 24

# Row 45
$ws.Cells.Item(45, 4).Value = 44676
$ws.Cells.Item(45, 11).Value = "Sutil De Gase"
$ws.Cells.Item(45, 12).Value = "Primera"
$ws.Cells.Item(45, 13).Value = 160
$ws.Cells.Item(45, 14).Value = 38000
$ws.Cells.Item(45, 15).Value = 40000
$ws.Cells.Item(45, 16).Value = 39000
$ws.Cells.Item(45, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(45, 18).Value = "Perú"
$ws.Cells.Item(45, 19).Value = 1625
$ws.Cells.Item(45, 20).Value = 24

# Row 46
$ws.Cells.Item(46, 4).Value = 44676
$ws.Cells.Item(46, 11).Value = "Tahití"
$ws.Cells.Item(46, 12).Value = "Primera"
$ws.Cells.Item(46, 13).Value = 200
$ws.Cells.Item(46, 14).Value = 14000
$ws.Cells.Item(46, 15).Value = 15000
$ws.Cells.Item(46, 16).Value = 14500
$ws.Cells.Item(46, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(46, 18).Value = "Perú"
$ws.Cells.Item(46, 19).Value = 806
$ws.Cells.Item(46, 20).Value = 18

# Row 47
$ws.Cells.Item(47, 4).Value = 44469
$ws.Cells.Item(47, 11).Value = "Sin especificar"
$ws.Cells.Item(47, 12).Value = "2a amarillo"
$ws.Cells.Item(47, 13).Value = 250
$ws.Cells.Item(47, 14).Value = 10000
$ws.Cells.Item(47, 15).Value = 11000
$ws.Cells.Item(47, 16).Value = 10500
$ws.Cells.Item(47, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(47, 18).Value = "Región Metropolitana"
$ws.Cells.Item(47, 19).Value = 525
$ws.Cells.Item(47, 20).Value = 20

# Row 48
$ws.Cells.Item(48, 4).Value = 44557
$ws.Cells.Item(48, 11).Value = "Tahití"
$ws.Cells.Item(48, 12).Value = "Primera"
$ws.Cells.Item(48, 13).Value = 200
$ws.Cells.Item(48, 14).Value = 40000
$ws.Cells.Item(48, 15).Value = 41000
$ws.Cells.Item(48, 16).Value = 40500
$ws.Cells.Item(48, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(48, 18).Value = "Perú"
$ws.Cells.Item(48, 19).Value = 1688
$ws.Cells.Item(48, 20).Value = 24

# Row 49
$ws.Cells.Item(49, 4).Value = 44230
$ws.Cells.Item(49, 11).Value = "Sin especificar"
$ws.Cells.Item(49, 12).Value = "2a plateado"
$ws.Cells.Item(49, 13).Value = 250
$ws.Cells.Item(49, 14).Value = 27000
$ws.Cells.Item(49, 15).Value = 28000
$ws.Cells.Item(49, 16).Value = 27500
$ws.Cells.Item(49, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(49, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(49, 19).Value = 1375
$ws.Cells.Item(49, 20).Value = 20

# Row 50
$ws.Cells.Item(50, 4).Value = 44372
$ws.Cells.Item(50, 11).Value = "Sutil De Gase"
$ws.Cells.Item(50, 12).Value = "Primera"
$ws.Cells.Item(50, 13).Value = 135
$ws.Cells.Item(50, 14).Value = 24000
$ws.Cells.Item(50, 15).Value = 25000
$ws.Cells.Item(50, 16).Value = 24370
$ws.Cells.Item(50, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(50, 18).Value = "Perú"
$ws.Cells.Item(50, 19).Value = 1015
$ws.Cells.Item(50, 20).Value = 24

# Row 51
$ws.Cells.Item(51, 4).Value = 44372
$ws.Cells.Item(51, 11).Value = "Tahití"
$ws.Cells.Item(51, 12).Value = "Primera"
$ws.Cells.Item(51, 13).Value = 150
$ws.Cells.Item(51, 14).Value = 25000
$ws.Cells.Item(51, 15).Value = 26000
$ws.Cells.Item(51, 16).Value = 25667
$ws.Cells.Item(51, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(51, 18).Value = "Perú"
$ws.Cells.Item(51, 19).Value = 1069
$ws.Cells.Item(51, 20).Value = 24

# Row 52
$ws.Cells.Item(52, 4).Value = 44172
$ws.Cells.Item(52, 11).Value = "Sutil De Gase"
$ws.Cells.Item(52, 12).Value = "Primera"
$ws.Cells.Item(52, 13).Value = 200
$ws.Cells.Item(52, 14).Value = 30000
$ws.Cells.Item(52, 15).Value = 31000
$ws.Cells.Item(52, 16).Value = 30500
$ws.Cells.Item(52, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(52, 18).Value = "Perú"
$ws.Cells.Item(52, 19).Value = 1271
$ws.Cells.Item(52, 20).Value = 24

# Row 53
$ws.Cells.Item(53, 4).Value = 44172
$ws.Cells.Item(53, 11).Value = "Tahití"
$ws.Cells.Item(53, 12).Value = "Primera"
$ws.Cells.Item(53, 13).Value = 360
$ws.Cells.Item(53, 14).Value = 24000
$ws.Cells.Item(53, 15).Value = 25000
$ws.Cells.Item(53, 16).Value = 24500
$ws.Cells.Item(53, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(53, 18).Value = "Perú"
$ws.Cells.Item(53, 19).Value = 1021
$ws.Cells.Item(53, 20).Value = 24

# Row 54
$ws.Cells.Item(54, 4).Value = 44209
$ws.Cells.Item(54, 11).Value = "Sin especificar"
$ws.Cells.Item(54, 12).Value = "1a amarillo"
$ws.Cells.Item(54, 13).Value = 300
$ws.Cells.Item(54, 14).Value = 29000
$ws.Cells.Item(54, 15).Value = 30000
$ws.Cells.Item(54, 16).Value = 29500
$ws.Cells.Item(54, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(54, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(54, 19).Value = 1475
$ws.Cells.Item(54, 20).Value = 20

# Row 55
$ws.Cells.Item(55, 4).Value = 44235
$ws.Cells.Item(55, 11).Value = "Sutil De Gase"
$ws.Cells.Item(55, 12).Value = "Primera"
$ws.Cells.Item(55, 13).Value = 250
$ws.Cells.Item(55, 14).Value = 21000
$ws.Cells.Item(55, 15).Value = 22000
$ws.Cells.Item(55, 16).Value = 21500
$ws.Cells.Item(55, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(55, 18).Value = "Perú"
$ws.Cells.Item(55, 19).Value = 896
$ws.Cells.Item(55, 20).Value = 24

# Row 56
$ws.Cells.Item(56, 4).Value = 44235
$ws.Cells.Item(56, 11).Value = "Tahití"
$ws.Cells.Item(56, 12).Value = "Primera"
$ws.Cells.Item(56, 13).Value = 300
$ws.Cells.Item(56, 14).Value = 22000
$ws.Cells.Item(56, 15).Value = 23000
$ws.Cells.Item(56, 16).Value = 22500
$ws.Cells.Item(56, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(56, 18).Value = "Perú"
$ws.Cells.Item(56, 19).Value = 938
$ws.Cells.Item(56, 20).Value = 24

# Row 57
$ws.Cells.Item(57, 4).Value = 44673
$ws.Cells.Item(57, 11).Value = "Sutil De Gase"
$ws.Cells.Item(57, 12).Value = "Primera"
$ws.Cells.Item(57, 13).Value = 300
$ws.Cells.Item(57, 14).Value = 38000
$ws.Cells.Item(57, 15).Value = 39000
$ws.Cells.Item(57, 16).Value = 38500
$ws.Cells.Item(57, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(57, 18).Value = "Perú"
$ws.Cells.Item(57, 19).Value = 1604
$ws.Cells.Item(57, 20).Value = 24

# Row 58
$ws.Cells.Item(58, 4).Value = 44636
$ws.Cells.Item(58, 11).Value = "Sin especificar"
$ws.Cells.Item(58, 12).Value = "2a amarillo"
$ws.Cells.Item(58, 13).Value = 300
$ws.Cells.Item(58, 14).Value = 32000
$ws.Cells.Item(58, 15).Value = 35000
$ws.Cells.Item(58, 16).Value = 33500
$ws.Cells.Item(58, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(58, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(58, 19).Value = 1675
$ws.Cells.Item(58, 20).Value = 20

# Row 59
$ws.Cells.Item(59, 4).Value = 44165
$ws.Cells.Item(59, 11).Value = "Sutil De Gase"
$ws.Cells.Item(59, 12).Value = "Primera"
$ws.Cells.Item(59, 13).Value = 200
$ws.Cells.Item(59, 14).Value = 31000
$ws.Cells.Item(59, 15).Value = 32000
$ws.Cells.Item(59, 16).Value = 31500
$ws.Cells.Item(59, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(59, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(59, 19).Value = 1312
$ws.Cells.Item(59, 20).Value = 24

# Row 60
$ws.Cells.Item(60, 4).Value = 44165
$ws.Cells.Item(60, 11).Value = "Tahití"
$ws.Cells.Item(60, 12).Value = "Primera"
$ws.Cells.Item(60, 13).Value = 340
$ws.Cells.Item(60, 14).Value = 28000
$ws.Cells.Item(60, 15).Value = 29000
$ws.Cells.Item(60, 16).Value = 28500
$ws.Cells.Item(60, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(60, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(60, 19).Value = 1188
$ws.Cells.Item(60, 20).Value = 24

# Row 61
$ws.Cells.Item(61, 4).Value = 44396
$ws.Cells.Item(61, 11).Value = "Sutil De Gase"
$ws.Cells.Item(61, 12).Value = "Primera"
$ws.Cells.Item(61, 13).Value = 200
$ws.Cells.Item(61, 14).Value = 30000
$ws.Cells.Item(61, 15).Value = 31000
$ws.Cells.Item(61, 16).Value = 30500
$ws.Cells.Item(61, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(61, 18).Value = "Perú"
$ws.Cells.Item(61, 19).Value = 1271
$ws.Cells.Item(61, 20).Value = 24

# Row 62
$ws.Cells.Item(62, 4).Value = 44342
$ws.Cells.Item(62, 11).Value = "Sin especificar"
$ws.Cells.Item(62, 12).Value = "2a plateado"
$ws.Cells.Item(62, 13).Value = 250
$ws.Cells.Item(62, 14).Value = 17000
$ws.Cells.Item(62, 15).Value = 18000
$ws.Cells.Item(62, 16).Value = 17500
$ws.Cells.Item(62, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(62, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(62, 19).Value = 875
$ws.Cells.Item(62, 20).Value = 20

# Row 63
$ws.Cells.Item(63, 4).Value = 44270
$ws.Cells.Item(63, 11).Value = "Tahití"
$ws.Cells.Item(63, 12).Value = "Primera"
$ws.Cells.Item(63, 13).Value = 250
$ws.Cells.Item(63, 14).Value = 29000
$ws.Cells.Item(63, 15).Value = 30000
$ws.Cells.Item(63, 16).Value = 29500
$ws.Cells.Item(63, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(63, 18).Value = "Perú"
$ws.Cells.Item(63, 19).Value = 1229
$ws.Cells.Item(63, 20).Value = 24

# Row 64
$ws.Cells.Item(64, 4).Value = 44340
$ws.Cells.Item(64, 11).Value = "Tahití"
$ws.Cells.Item(64, 12).Value = "Primera"
$ws.Cells.Item(64, 13).Value = 250
$ws.Cells.Item(64, 14).Value = 25000
$ws.Cells.Item(64, 15).Value = 26000
$ws.Cells.Item(64, 16).Value = 25500
$ws.Cells.Item(64, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(64, 18).Value = "Perú"
$ws.Cells.Item(64, 19).Value = 1062
$ws.Cells.Item(64, 20).Value = 24

# Row 65
$ws.Cells.Item(65, 4).Value = 44552
$ws.Cells.Item(65, 11).Value = "Sin especificar"
$ws.Cells.Item(65, 12).Value = "2a amarillo"
$ws.Cells.Item(65, 13).Value = 300
$ws.Cells.Item(65, 14).Value = 21000
$ws.Cells.Item(65, 15).Value = 22000
$ws.Cells.Item(65, 16).Value = 21500
$ws.Cells.Item(65, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(65, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(65, 19).Value = 1075
$ws.Cells.Item(65, 20).Value = 20

# Row 66
$ws.Cells.Item(66, 4).Value = 44298
$ws.Cells.Item(66, 11).Value = "Sutil De Gase"
$ws.Cells.Item(66, 12).Value = "Primera"
$ws.Cells.Item(66, 13).Value = 160
$ws.Cells.Item(66, 14).Value = 31000
$ws.Cells.Item(66, 15).Value = 32000
$ws.Cells.Item(66, 16).Value = 31500
$ws.Cells.Item(66, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(66, 18).Value = "Perú"
$ws.Cells.Item(66, 19).Value = 1312
$ws.Cells.Item(66, 20).Value = 24

# Row 67
$ws.Cells.Item(67, 4).Value = 44298
$ws.Cells.Item(67, 11).Value = "Tahití"
$ws.Cells.Item(67, 12).Value = "Primera"
$ws.Cells.Item(67, 13).Value = 300
$ws.Cells.Item(67, 14).Value = 27000
$ws.Cells.Item(67, 15).Value = 28000
$ws.Cells.Item(67, 16).Value = 27500
$ws.Cells.Item(67, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(67, 18).Value = "Perú"
$ws.Cells.Item(67, 19).Value = 1146
$ws.Cells.Item(67, 20).Value = 24

# Row 68
$ws.Cells.Item(68, 4).Value = 44277
$ws.Cells.Item(68, 11).Value = "Tahití"
$ws.Cells.Item(68, 12).Value = "Primera"
$ws.Cells.Item(68, 13).Value = 300
$ws.Cells.Item(68, 14).Value = 29000
$ws.Cells.Item(68, 15).Value = 30000
$ws.Cells.Item(68, 16).Value = 29500
$ws.Cells.Item(68, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(68, 18).Value = "Perú"
$ws.Cells.Item(68, 19).Value = 1229
$ws.Cells.Item(68, 20).Value = 24

# Row 69
$ws.Cells.Item(69, 4).Value = 44399
$ws.Cells.Item(69, 11).Value = "Sin especificar"
$ws.Cells.Item(69, 12).Value = "2a amarillo"
$ws.Cells.Item(69, 13).Value = 240
$ws.Cells.Item(69, 14).Value = 11000
$ws.Cells.Item(69, 15).Value = 12000
$ws.Cells.Item(69, 16).Value = 11500
$ws.Cells.Item(69, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(69, 18).Value = "Región Metropolitana"
$ws.Cells.Item(69, 19).Value = 575
$ws.Cells.Item(69, 20).Value = 20

# Row 70
$ws.Cells.Item(70, 4).Value = 44329
$ws.Cells.Item(70, 11).Value = "Sin especificar"
$ws.Cells.Item(70, 12).Value = "2a amarillo"
$ws.Cells.Item(70, 13).Value = 250
$ws.Cells.Item(70, 14).Value = 25000
$ws.Cells.Item(70, 15).Value = 26000
$ws.Cells.Item(70, 16).Value = 25500
$ws.Cells.Item(70, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(70, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(70, 19).Value = 1275
$ws.Cells.Item(70, 20).Value = 20

# Row 71
$ws.Cells.Item(71, 4).Value = 44420
$ws.Cells.Item(71, 11).Value = "Sin especificar"
$ws.Cells.Item(71, 12).Value = "3a amarillo"
$ws.Cells.Item(71, 13).Value = 250
$ws.Cells.Item(71, 14).Value = 10000
$ws.Cells.Item(71, 15).Value = 11000
$ws.Cells.Item(71, 16).Value = 10500
$ws.Cells.Item(71, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(71, 18).Value = "Región Metropolitana"
$ws.Cells.Item(71, 19).Value = 525
$ws.Cells.Item(71, 20).Value = 20

# Row 72
$ws.Cells.Item(72, 4).Value = 44384
$ws.Cells.Item(72, 11).Value = "Sin especificar"
$ws.Cells.Item(72, 12).Value = "2a amarillo"
$ws.Cells.Item(72, 13).Value = 150
$ws.Cells.Item(72, 14).Value = 11000
$ws.Cells.Item(72, 15).Value = 12000
$ws.Cells.Item(72, 16).Value = 11667
$ws.Cells.Item(72, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(72, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(72, 19).Value = 583
$ws.Cells.Item(72, 20).Value = 20

# Row 73
$ws.Cells.Item(73, 4).Value = 44384
$ws.Cells.Item(73, 11).Value = "Sutil De Gase"
$ws.Cells.Item(73, 12).Value = "Primera"
$ws.Cells.Item(73, 13).Value = 140
$ws.Cells.Item(73, 14).Value = 32000
$ws.Cells.Item(73, 15).Value = 33000
$ws.Cells.Item(73, 16).Value = 32429
$ws.Cells.Item(73, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(73, 18).Value = "Perú"
$ws.Cells.Item(73, 19).Value = 1351
$ws.Cells.Item(73, 20).Value = 24

# Row 74
$ws.Cells.Item(74, 4).Value = 44214
$ws.Cells.Item(74, 11).Value = "Sutil De Gase"
$ws.Cells.Item(74, 12).Value = "Primera"
$ws.Cells.Item(74, 13).Value = 250
$ws.Cells.Item(74, 14).Value = 29000
$ws.Cells.Item(74, 15).Value = 30000
$ws.Cells.Item(74, 16).Value = 29500
$ws.Cells.Item(74, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(74, 18).Value = "Perú"
$ws.Cells.Item(74, 19).Value = 1229
$ws.Cells.Item(74, 20).Value = 24

# Row 75
$ws.Cells.Item(75, 4).Value = 44214
$ws.Cells.Item(75, 11).Value = "Tahití"
$ws.Cells.Item(75, 12).Value = "Primera"
$ws.Cells.Item(75, 13).Value = 200
$ws.Cells.Item(75, 14).Value = 26000
$ws.Cells.Item(75, 15).Value = 27000
$ws.Cells.Item(75, 16).Value = 26500
$ws.Cells.Item(75, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(75, 18).Value = "Perú"
$ws.Cells.Item(75, 19).Value = 1104
$ws.Cells.Item(75, 20).Value = 24

# Row 76
$ws.Cells.Item(76, 4).Value = 44319
$ws.Cells.Item(76, 11).Value = "Sutil De Gase"
$ws.Cells.Item(76, 12).Value = "Primera"
$ws.Cells.Item(76, 13).Value = 250
$ws.Cells.Item(76, 14).Value = 27000
$ws.Cells.Item(76, 15).Value = 28000
$ws.Cells.Item(76, 16).Value = 27500
$ws.Cells.Item(76, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(76, 18).Value = "Perú"
$ws.Cells.Item(76, 19).Value = 1146
$ws.Cells.Item(76, 20).Value = 24

# Row 77
$ws.Cells.Item(77, 4).Value = 44319
$ws.Cells.Item(77, 11).Value = "Tahití"
$ws.Cells.Item(77, 12).Value = "Primera"
$ws.Cells.Item(77, 13).Value = 300
$ws.Cells.Item(77, 14).Value = 23000
$ws.Cells.Item(77, 15).Value = 24000
$ws.Cells.Item(77, 16).Value = 23500
$ws.Cells.Item(77, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(77, 18).Value = "Perú"
$ws.Cells.Item(77, 19).Value = 979
$ws.Cells.Item(77, 20).Value = 24

# Row 78
$ws.Cells.Item(78, 4).Value = 44242
$ws.Cells.Item(78, 11).Value = "Sutil De Gase"
$ws.Cells.Item(78, 12).Value = "Primera"
$ws.Cells.Item(78, 13).Value = 250
$ws.Cells.Item(78, 14).Value = 22000
$ws.Cells.Item(78, 15).Value = 23000
$ws.Cells.Item(78, 16).Value = 22500
$ws.Cells.Item(78, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(78, 18).Value = "Perú"
$ws.Cells.Item(78, 19).Value = 938
$ws.Cells.Item(78, 20).Value = 24

# Row 79
$ws.Cells.Item(79, 4).Value = 44242
$ws.Cells.Item(79, 11).Value = "Tahití"
$ws.Cells.Item(79, 12).Value = "Primera"
$ws.Cells.Item(79, 13).Value = 300
$ws.Cells.Item(79, 14).Value = 21000
$ws.Cells.Item(79, 15).Value = 22000
$ws.Cells.Item(79, 16).Value = 21500
$ws.Cells.Item(79, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(79, 18).Value = "Perú"
$ws.Cells.Item(79, 19).Value = 896
$ws.Cells.Item(79, 20).Value = 24

# Row 80
$ws.Cells.Item(80, 4).Value = 44265
$ws.Cells.Item(80, 11).Value = "Sin especificar"
$ws.Cells.Item(80, 12).Value = "2a amarillo"
$ws.Cells.Item(80, 13).Value = 250
$ws.Cells.Item(80, 14).Value = 26000
$ws.Cells.Item(80, 15).Value = 27000
$ws.Cells.Item(80, 16).Value = 26500
$ws.Cells.Item(80, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(80, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(80, 19).Value = 1325
$ws.Cells.Item(80, 20).Value = 20

# Row 81
$ws.Cells.Item(81, 4).Value = 44657
$ws.Cells.Item(81, 11).Value = "Tahití"
$ws.Cells.Item(81, 12).Value = "Primera"
$ws.Cells.Item(81, 13).Value = 300
$ws.Cells.Item(81, 14).Value = 21000
$ws.Cells.Item(81, 15).Value = 22000
$ws.Cells.Item(81, 16).Value = 21500
$ws.Cells.Item(81, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(81, 18).Value = "Perú"
$ws.Cells.Item(81, 19).Value = 1194
$ws.Cells.Item(81, 20).Value = 18

# Row 82
$ws.Cells.Item(82, 4).Value = 44599
$ws.Cells.Item(82, 11).Value = "Tahití"
$ws.Cells.Item(82, 12).Value = "Primera"
$ws.Cells.Item(82, 13).Value = 300
$ws.Cells.Item(82, 14).Value = 35000
$ws.Cells.Item(82, 15).Value = 36000
$ws.Cells.Item(82, 16).Value = 35500
$ws.Cells.Item(82, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(82, 18).Value = "Perú"
$ws.Cells.Item(82, 19).Value = 1479
$ws.Cells.Item(82, 20).Value = 24

# Row 83
$ws.Cells.Item(83, 4).Value = 44344
$ws.Cells.Item(83, 11).Value = "Sutil De Gase"
$ws.Cells.Item(83, 12).Value = "Primera"
$ws.Cells.Item(83, 13).Value = 250
$ws.Cells.Item(83, 14).Value = 25000
$ws.Cells.Item(83, 15).Value = 26000
$ws.Cells.Item(83, 16).Value = 25500
$ws.Cells.Item(83, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(83, 18).Value = "Perú"
$ws.Cells.Item(83, 19).Value = 1062
$ws.Cells.Item(83, 20).Value = 24

# Row 84
$ws.Cells.Item(84, 4).Value = 44341
$ws.Cells.Item(84, 11).Value = "Sin especificar"
$ws.Cells.Item(84, 12).Value = "3a plateado"
$ws.Cells.Item(84, 13).Value = 250
$ws.Cells.Item(84, 14).Value = 10000
$ws.Cells.Item(84, 15).Value = 11000
$ws.Cells.Item(84, 16).Value = 10500
$ws.Cells.Item(84, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(84, 18).Value = "Región Metropolitana"
$ws.Cells.Item(84, 19).Value = 525
$ws.Cells.Item(84, 20).Value = 20

# Row 85
$ws.Cells.Item(85, 4).Value = 44356
$ws.Cells.Item(85, 11).Value = "Sin especificar"
$ws.Cells.Item(85, 12).Value = "2a amarillo"
$ws.Cells.Item(85, 13).Value = 200
$ws.Cells.Item(85, 14).Value = 10000
$ws.Cells.Item(85, 15).Value = 11000
$ws.Cells.Item(85, 16).Value = 10500
$ws.Cells.Item(85, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(85, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(85, 19).Value = 525
$ws.Cells.Item(85, 20).Value = 20

# Row 86
$ws.Cells.Item(86, 4).Value = 44356
$ws.Cells.Item(86, 11).Value = "Sutil De Gase"
$ws.Cells.Item(86, 12).Value = "Primera"
$ws.Cells.Item(86, 13).Value = 200
$ws.Cells.Item(86, 14).Value = 24000
$ws.Cells.Item(86, 15).Value = 25000
$ws.Cells.Item(86, 16).Value = 24500
$ws.Cells.Item(86, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(86, 18).Value = "Perú"
$ws.Cells.Item(86, 19).Value = 1021
$ws.Cells.Item(86, 20).Value = 24

# Row 87
$ws.Cells.Item(87, 4).Value = 44356
$ws.Cells.Item(87, 11).Value = "Tahití"
$ws.Cells.Item(87, 12).Value = "Primera"
$ws.Cells.Item(87, 13).Value = 250
$ws.Cells.Item(87, 14).Value = 24000
$ws.Cells.Item(87, 15).Value = 25000
$ws.Cells.Item(87, 16).Value = 24500
$ws.Cells.Item(87, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(87, 18).Value = "Perú"
$ws.Cells.Item(87, 19).Value = 1021
$ws.Cells.Item(87, 20).Value = 24

# Row 88
$ws.Cells.Item(88, 4).Value = 44494
$ws.Cells.Item(88, 11).Value = "Sutil De Gase"
$ws.Cells.Item(88, 12).Value = "Primera"
$ws.Cells.Item(88, 13).Value = 160
$ws.Cells.Item(88, 14).Value = 54000
$ws.Cells.Item(88, 15).Value = 55000
$ws.Cells.Item(88, 16).Value = 54500
$ws.Cells.Item(88, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(88, 18).Value = "Perú"
$ws.Cells.Item(88, 19).Value = 2271
$ws.Cells.Item(88, 20).Value = 24

# Row 89
$ws.Cells.Item(89, 4).Value = 44629
$ws.Cells.Item(89, 11).Value = "Sin especificar"
$ws.Cells.Item(89, 12).Value = "2a amarillo"
$ws.Cells.Item(89, 13).Value = 270
$ws.Cells.Item(89, 14).Value = 29000
$ws.Cells.Item(89, 15).Value = 30000
$ws.Cells.Item(89, 16).Value = 29500
$ws.Cells.Item(89, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(89, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(89, 19).Value = 1475
$ws.Cells.Item(89, 20).Value = 20

# Row 90
$ws.Cells.Item(90, 4).Value = 44427
$ws.Cells.Item(90, 11).Value = "Sin especificar"
$ws.Cells.Item(90, 12).Value = "2a amarillo"
$ws.Cells.Item(90, 13).Value = 250
$ws.Cells.Item(90, 14).Value = 10000
$ws.Cells.Item(90, 15).Value = 11000
$ws.Cells.Item(90, 16).Value = 10500
$ws.Cells.Item(90, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(90, 18).Value = "Región Metropolitana"
$ws.Cells.Item(90, 19).Value = 525
$ws.Cells.Item(90, 20).Value = 20

# Row 91
$ws.Cells.Item(91, 4).Value = 44573
$ws.Cells.Item(91, 11).Value = "Sin especificar"
$ws.Cells.Item(91, 12).Value = "1a amarillo"
$ws.Cells.Item(91, 13).Value = 270
$ws.Cells.Item(91, 14).Value = 27000
$ws.Cells.Item(91, 15).Value = 28000
$ws.Cells.Item(91, 16).Value = 27500
$ws.Cells.Item(91, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(91, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(91, 19).Value = 1375
$ws.Cells.Item(91, 20).Value = 20

# Row 92
$ws.Cells.Item(92, 4).Value = 44487
$ws.Cells.Item(92, 11).Value = "Tahití"
$ws.Cells.Item(92, 12).Value = "Primera"
$ws.Cells.Item(92, 13).Value = 130
$ws.Cells.Item(92, 14).Value = 47000
$ws.Cells.Item(92, 15).Value = 48000
$ws.Cells.Item(92, 16).Value = 47500
$ws.Cells.Item(92, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(92, 18).Value = "Perú"
$ws.Cells.Item(92, 19).Value = 1979
$ws.Cells.Item(92, 20).Value = 24

# Row 93
$ws.Cells.Item(93, 4).Value = 44412
$ws.Cells.Item(93, 11).Value = "Sin especificar"
$ws.Cells.Item(93, 12).Value = "2a amarillo"
$ws.Cells.Item(93, 13).Value = 300
$ws.Cells.Item(93, 14).Value = 9000
$ws.Cells.Item(93, 15).Value = 10000
$ws.Cells.Item(93, 16).Value = 9500
$ws.Cells.Item(93, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(93, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(93, 19).Value = 475
$ws.Cells.Item(93, 20).Value = 20

# Row 94
$ws.Cells.Item(94, 4).Value = 44204
$ws.Cells.Item(94, 11).Value = "Sutil De Gase"
$ws.Cells.Item(94, 12).Value = "Primera"
$ws.Cells.Item(94, 13).Value = 250
$ws.Cells.Item(94, 14).Value = 29000
$ws.Cells.Item(94, 15).Value = 30000
$ws.Cells.Item(94, 16).Value = 29500
$ws.Cells.Item(94, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(94, 18).Value = "Perú"
$ws.Cells.Item(94, 19).Value = 1229
$ws.Cells.Item(94, 20).Value = 24

# Row 95
$ws.Cells.Item(95, 4).Value = 44204
$ws.Cells.Item(95, 11).Value = "Tahití"
$ws.Cells.Item(95, 12).Value = "Primera"
$ws.Cells.Item(95, 13).Value = 300
$ws.Cells.Item(95, 14).Value = 27000
$ws.Cells.Item(95, 15).Value = 28000
$ws.Cells.Item(95, 16).Value = 27500
$ws.Cells.Item(95, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(95, 18).Value = "Perú"
$ws.Cells.Item(95, 19).Value = 1146
$ws.Cells.Item(95, 20).Value = 24

# Row 96
$ws.Cells.Item(96, 4).Value = 44300
$ws.Cells.Item(96, 11).Value = "Sin especificar"
$ws.Cells.Item(96, 12).Value = "2a amarillo"
$ws.Cells.Item(96, 13).Value = 270
$ws.Cells.Item(96, 14).Value = 20000
$ws.Cells.Item(96, 15).Value = 21000
$ws.Cells.Item(96, 16).Value = 20500
$ws.Cells.Item(96, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(96, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(96, 19).Value = 1025
$ws.Cells.Item(96, 20).Value = 20

# Row 97
$ws.Cells.Item(97, 4).Value = 44566
$ws.Cells.Item(97, 11).Value = "Sin especificar"
$ws.Cells.Item(97, 12).Value = "2a amarillo"
$ws.Cells.Item(97, 13).Value = 300
$ws.Cells.Item(97, 14).Value = 24000
$ws.Cells.Item(97, 15).Value = 25000
$ws.Cells.Item(97, 16).Value = 24500
$ws.Cells.Item(97, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(97, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(97, 19).Value = 1225
$ws.Cells.Item(97, 20).Value = 20

# Row 98
$ws.Cells.Item(98, 4).Value = 44392
$ws.Cells.Item(98, 11).Value = "Sin especificar"
$ws.Cells.Item(98, 12).Value = "2a amarillo"
$ws.Cells.Item(98, 13).Value = 270
$ws.Cells.Item(98, 14).Value = 11000
$ws.Cells.Item(98, 15).Value = 12000
$ws.Cells.Item(98, 16).Value = 11500
$ws.Cells.Item(98, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(98, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(98, 19).Value = 575
$ws.Cells.Item(98, 20).Value = 20

# Row 99
$ws.Cells.Item(99, 4).Value = 44308
$ws.Cells.Item(99, 11).Value = "Sin especificar"
$ws.Cells.Item(99, 12).Value = "2a plateado"
$ws.Cells.Item(99, 13).Value = 250
$ws.Cells.Item(99, 14).Value = 19000
$ws.Cells.Item(99, 15).Value = 20000
$ws.Cells.Item(99, 16).Value = 19500
$ws.Cells.Item(99, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(99, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(99, 19).Value = 975
$ws.Cells.Item(99, 20).Value = 20

# Row 100
$ws.Cells.Item(100, 4).Value = 44330
$ws.Cells.Item(100, 11).Value = "Tahití"
$ws.Cells.Item(100, 12).Value = "Primera"
$ws.Cells.Item(100, 13).Value = 250
$ws.Cells.Item(100, 14).Value = 27000
$ws.Cells.Item(100, 15).Value = 28000
$ws.Cells.Item(100, 16).Value = 27500
$ws.Cells.Item(100, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(100, 18).Value = "Perú"
$ws.Cells.Item(100, 19).Value = 1146
$ws.Cells.Item(100, 20).Value = 24

# Row 101
$ws.Cells.Item(101, 4).Value = 44200
$ws.Cells.Item(101, 11).Value = "Sutil De Gase"
$ws.Cells.Item(101, 12).Value = "Primera"
$ws.Cells.Item(101, 13).Value = 270
$ws.Cells.Item(101, 14).Value = 25000
$ws.Cells.Item(101, 15).Value = 26000
$ws.Cells.Item(101, 16).Value = 25500
$ws.Cells.Item(101, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(101, 18).Value = "Perú"
$ws.Cells.Item(101, 19).Value = 1062
$ws.Cells.Item(101, 20).Value = 24

# Row 102
$ws.Cells.Item(102, 4).Value = 44200
$ws.Cells.Item(102, 11).Value = "Tahití"
$ws.Cells.Item(102, 12).Value = "Primera"
$ws.Cells.Item(102, 13).Value = 300
$ws.Cells.Item(102, 14).Value = 29000
$ws.Cells.Item(102, 15).Value = 30000
$ws.Cells.Item(102, 16).Value = 29500
$ws.Cells.Item(102, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(102, 18).Value = "Perú"
$ws.Cells.Item(102, 19).Value = 1229
$ws.Cells.Item(102, 20).Value = 24

# Row 103
$ws.Cells.Item(103, 4).Value = 44610
$ws.Cells.Item(103, 11).Value = "Tahití"
$ws.Cells.Item(103, 12).Value = "Primera"
$ws.Cells.Item(103, 13).Value = 300
$ws.Cells.Item(103, 14).Value = 35000
$ws.Cells.Item(103, 15).Value = 36000
$ws.Cells.Item(103, 16).Value = 35500
$ws.Cells.Item(103, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(103, 18).Value = "Perú"
$ws.Cells.Item(103, 19).Value = 1479
$ws.Cells.Item(103, 20).Value = 24

# Row 104
$ws.Cells.Item(104, 4).Value = 44536
$ws.Cells.Item(104, 11).Value = "Sutil De Gase"
$ws.Cells.Item(104, 12).Value = "Segunda"
$ws.Cells.Item(104, 13).Value = 200
$ws.Cells.Item(104, 14).Value = 19000
$ws.Cells.Item(104, 15).Value = 20000
$ws.Cells.Item(104, 16).Value = 19500
$ws.Cells.Item(104, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(104, 18).Value = "Perú"
$ws.Cells.Item(104, 19).Value = 812
$ws.Cells.Item(104, 20).Value = 24

# Row 105
$ws.Cells.Item(105, 4).Value = 44536
$ws.Cells.Item(105, 11).Value = "Tahití"
$ws.Cells.Item(105, 12).Value = "Primera"
$ws.Cells.Item(105, 13).Value = 250
$ws.Cells.Item(105, 14).Value = 29000
$ws.Cells.Item(105, 15).Value = 30000
$ws.Cells.Item(105, 16).Value = 29500
$ws.Cells.Item(105, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(105, 18).Value = "Perú"
$ws.Cells.Item(105, 19).Value = 1229
$ws.Cells.Item(105, 20).Value = 24

# Row 106
$ws.Cells.Item(106, 4).Value = 44169
$ws.Cells.Item(106, 11).Value = "Sutil De Gase"
$ws.Cells.Item(106, 12).Value = "Primera"
$ws.Cells.Item(106, 13).Value = 270
$ws.Cells.Item(106, 14).Value = 29000
$ws.Cells.Item(106, 15).Value = 31000
$ws.Cells.Item(106, 16).Value = 30000
$ws.Cells.Item(106, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(106, 18).Value = "Perú"
$ws.Cells.Item(106, 19).Value = 1250
$ws.Cells.Item(106, 20).Value = 24

# Row 107
$ws.Cells.Item(107, 4).Value = 44169
$ws.Cells.Item(107, 11).Value = "Tahití"
$ws.Cells.Item(107, 12).Value = "Primera"
$ws.Cells.Item(107, 13).Value = 360
$ws.Cells.Item(107, 14).Value = 23000
$ws.Cells.Item(107, 15).Value = 24000
$ws.Cells.Item(107, 16).Value = 23500
$ws.Cells.Item(107, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(107, 18).Value = "Perú"
$ws.Cells.Item(107, 19).Value = 979
$ws.Cells.Item(107, 20).Value = 24

# Row 108
$ws.Cells.Item(108, 4).Value = 44309
$ws.Cells.Item(108, 11).Value = "Sutil De Gase"
$ws.Cells.Item(108, 12).Value = "Primera"
$ws.Cells.Item(108, 13).Value = 250
$ws.Cells.Item(108, 14).Value = 31000
$ws.Cells.Item(108, 15).Value = 32000
$ws.Cells.Item(108, 16).Value = 31500
$ws.Cells.Item(108, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(108, 18).Value = "Perú"
$ws.Cells.Item(108, 19).Value = 1312
$ws.Cells.Item(108, 20).Value = 24

# Row 109
$ws.Cells.Item(109, 4).Value = 44309
$ws.Cells.Item(109, 11).Value = "Tahití"
$ws.Cells.Item(109, 12).Value = "Primera"
$ws.Cells.Item(109, 13).Value = 300
$ws.Cells.Item(109, 14).Value = 27000
$ws.Cells.Item(109, 15).Value = 28000
$ws.Cells.Item(109, 16).Value = 27500
$ws.Cells.Item(109, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(109, 18).Value = "Perú"
$ws.Cells.Item(109, 19).Value = 1146
$ws.Cells.Item(109, 20).Value = 24

# Row 110
$ws.Cells.Item(110, 4).Value = 44371
$ws.Cells.Item(110, 11).Value = "Sin especificar"
$ws.Cells.Item(110, 12).Value = "2a amarillo"
$ws.Cells.Item(110, 13).Value = 200
$ws.Cells.Item(110, 14).Value = 11000
$ws.Cells.Item(110, 15).Value = 12000
$ws.Cells.Item(110, 16).Value = 11600
$ws.Cells.Item(110, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(110, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(110, 19).Value = 580
$ws.Cells.Item(110, 20).Value = 20

# Row 111
$ws.Cells.Item(111, 4).Value = 44237
$ws.Cells.Item(111, 11).Value = "Sin especificar"
$ws.Cells.Item(111, 12).Value = "2a amarillo"
$ws.Cells.Item(111, 13).Value = 250
$ws.Cells.Item(111, 14).Value = 26000
$ws.Cells.Item(111, 15).Value = 27000
$ws.Cells.Item(111, 16).Value = 26500
$ws.Cells.Item(111, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(111, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(111, 19).Value = 1325
$ws.Cells.Item(111, 20).Value = 20

# Row 112
$ws.Cells.Item(112, 4).Value = 44237
$ws.Cells.Item(112, 11).Value = "Sutil De Gase"
$ws.Cells.Item(112, 12).Value = "Primera"
$ws.Cells.Item(112, 13).Value = 200
$ws.Cells.Item(112, 14).Value = 21000
$ws.Cells.Item(112, 15).Value = 22000
$ws.Cells.Item(112, 16).Value = 21500
$ws.Cells.Item(112, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(112, 18).Value = "Perú"
$ws.Cells.Item(112, 19).Value = 896
$ws.Cells.Item(112, 20).Value = 24

# Row 113
$ws.Cells.Item(113, 4).Value = 44237
$ws.Cells.Item(113, 11).Value = "Tahití"
$ws.Cells.Item(113, 12).Value = "Primera"
$ws.Cells.Item(113, 13).Value = 200
$ws.Cells.Item(113, 14).Value = 22000
$ws.Cells.Item(113, 15).Value = 23000
$ws.Cells.Item(113, 16).Value = 22500
$ws.Cells.Item(113, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(113, 18).Value = "Perú"
$ws.Cells.Item(113, 19).Value = 938
$ws.Cells.Item(113, 20).Value = 24

# Row 114
$ws.Cells.Item(114, 4).Value = 44333
$ws.Cells.Item(114, 11).Value = "Tahití"
$ws.Cells.Item(114, 12).Value = "Primera"
$ws.Cells.Item(114, 13).Value = 250
$ws.Cells.Item(114, 14).Value = 27000
$ws.Cells.Item(114, 15).Value = 28000
$ws.Cells.Item(114, 16).Value = 27500
$ws.Cells.Item(114, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(114, 18).Value = "Perú"
$ws.Cells.Item(114, 19).Value = 1146
$ws.Cells.Item(114, 20).Value = 24

# Row 115
$ws.Cells.Item(115, 4).Value = 44417
$ws.Cells.Item(115, 11).Value = "Sutil De Gase"
$ws.Cells.Item(115, 12).Value = "Primera"
$ws.Cells.Item(115, 13).Value = 300
$ws.Cells.Item(115, 14).Value = 32000
$ws.Cells.Item(115, 15).Value = 33000
$ws.Cells.Item(115, 16).Value = 32500
$ws.Cells.Item(115, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(115, 18).Value = "Perú"
$ws.Cells.Item(115, 19).Value = 1354
$ws.Cells.Item(115, 20).Value = 24

# Row 116
$ws.Cells.Item(116, 4).Value = 44417
$ws.Cells.Item(116, 11).Value = "Tahití"
$ws.Cells.Item(116, 12).Value = "Primera"
$ws.Cells.Item(116, 13).Value = 300
$ws.Cells.Item(116, 14).Value = 30000
$ws.Cells.Item(116, 15).Value = 31000
$ws.Cells.Item(116, 16).Value = 30500
$ws.Cells.Item(116, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(116, 18).Value = "Perú"
$ws.Cells.Item(116, 19).Value = 1271
$ws.Cells.Item(116, 20).Value = 24

# Row 117
$ws.Cells.Item(117, 4).Value = 44382
$ws.Cells.Item(117, 11).Value = "Sutil De Gase"
$ws.Cells.Item(117, 12).Value = "Primera"
$ws.Cells.Item(117, 13).Value = 180
$ws.Cells.Item(117, 14).Value = 32000
$ws.Cells.Item(117, 15).Value = 33000
$ws.Cells.Item(117, 16).Value = 32556
$ws.Cells.Item(117, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(117, 18).Value = "Perú"
$ws.Cells.Item(117, 19).Value = 1356
$ws.Cells.Item(117, 20).Value = 24

# Row 118
$ws.Cells.Item(118, 4).Value = 44260
$ws.Cells.Item(118, 11).Value = "Sutil De Gase"
$ws.Cells.Item(118, 12).Value = "Primera"
$ws.Cells.Item(118, 13).Value = 200
$ws.Cells.Item(118, 14).Value = 27000
$ws.Cells.Item(118, 15).Value = 28000
$ws.Cells.Item(118, 16).Value = 27500
$ws.Cells.Item(118, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(118, 18).Value = "Perú"
$ws.Cells.Item(118, 19).Value = 1146
$ws.Cells.Item(118, 20).Value = 24

# Row 119
$ws.Cells.Item(119, 4).Value = 44260
$ws.Cells.Item(119, 11).Value = "Tahití"
$ws.Cells.Item(119, 12).Value = "Primera"
$ws.Cells.Item(119, 13).Value = 300
$ws.Cells.Item(119, 14).Value = 22000
$ws.Cells.Item(119, 15).Value = 23000
$ws.Cells.Item(119, 16).Value = 22500
$ws.Cells.Item(119, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(119, 18).Value = "Perú"
$ws.Cells.Item(119, 19).Value = 938
$ws.Cells.Item(119, 20).Value = 24

# Row 120
$ws.Cells.Item(120, 4).Value = 44588
$ws.Cells.Item(120, 11).Value = "Sin especificar"
$ws.Cells.Item(120, 12).Value = "2a amarillo"
$ws.Cells.Item(120, 13).Value = 300
$ws.Cells.Item(120, 14).Value = 24000
$ws.Cells.Item(120, 15).Value = 25000
$ws.Cells.Item(120, 16).Value = 24500
$ws.Cells.Item(120, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(120, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(120, 19).Value = 1225
$ws.Cells.Item(120, 20).Value = 20

# Row 121
$ws.Cells.Item(121, 4).Value = 44179
$ws.Cells.Item(121, 11).Value = "Sutil De Gase"
$ws.Cells.Item(121, 12).Value = "Primera"
$ws.Cells.Item(121, 13).Value = 250
$ws.Cells.Item(121, 14).Value = 30000
$ws.Cells.Item(121, 15).Value = 31000
$ws.Cells.Item(121, 16).Value = 30500
$ws.Cells.Item(121, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(121, 18).Value = "Perú"
$ws.Cells.Item(121, 19).Value = 1271
$ws.Cells.Item(121, 20).Value = 24

# Row 122
$ws.Cells.Item(122, 4).Value = 44179
$ws.Cells.Item(122, 11).Value = "Tahití"
$ws.Cells.Item(122, 12).Value = "Primera"
$ws.Cells.Item(122, 13).Value = 360
$ws.Cells.Item(122, 14).Value = 24000
$ws.Cells.Item(122, 15).Value = 25000
$ws.Cells.Item(122, 16).Value = 24500
$ws.Cells.Item(122, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(122, 18).Value = "Perú"
$ws.Cells.Item(122, 19).Value = 1021
$ws.Cells.Item(122, 20).Value = 24

# Row 123
$ws.Cells.Item(123, 4).Value = 44522
$ws.Cells.Item(123, 11).Value = "Sutil De Gase"
$ws.Cells.Item(123, 12).Value = "Primera"
$ws.Cells.Item(123, 13).Value = 200
$ws.Cells.Item(123, 14).Value = 20000
$ws.Cells.Item(123, 15).Value = 21000
$ws.Cells.Item(123, 16).Value = 20500
$ws.Cells.Item(123, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(123, 18).Value = "Perú"
$ws.Cells.Item(123, 19).Value = 854
$ws.Cells.Item(123, 20).Value = 24

# Row 124
$ws.Cells.Item(124, 4).Value = 44522
$ws.Cells.Item(124, 11).Value = "Tahití"
$ws.Cells.Item(124, 12).Value = "Primera"
$ws.Cells.Item(124, 13).Value = 200
$ws.Cells.Item(124, 14).Value = 24000
$ws.Cells.Item(124, 15).Value = 25000
$ws.Cells.Item(124, 16).Value = 24500
$ws.Cells.Item(124, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(124, 18).Value = "Perú"
$ws.Cells.Item(124, 19).Value = 1021
$ws.Cells.Item(124, 20).Value = 24

# Row 125
$ws.Cells.Item(125, 4).Value = 44225
$ws.Cells.Item(125, 11).Value = "Sutil De Gase"
$ws.Cells.Item(125, 12).Value = "Primera"
$ws.Cells.Item(125, 13).Value = 250
$ws.Cells.Item(125, 14).Value = 25000
$ws.Cells.Item(125, 15).Value = 26000
$ws.Cells.Item(125, 16).Value = 25500
$ws.Cells.Item(125, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(125, 18).Value = "Perú"
$ws.Cells.Item(125, 19).Value = 1062
$ws.Cells.Item(125, 20).Value = 24

# Row 126
$ws.Cells.Item(126, 4).Value = 44225
$ws.Cells.Item(126, 11).Value = "Tahití"
$ws.Cells.Item(126, 12).Value = "Primera"
$ws.Cells.Item(126, 13).Value = 250
$ws.Cells.Item(126, 14).Value = 25000
$ws.Cells.Item(126, 15).Value = 26000
$ws.Cells.Item(126, 16).Value = 25500
$ws.Cells.Item(126, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(126, 18).Value = "Perú"
$ws.Cells.Item(126, 19).Value = 1062
$ws.Cells.Item(126, 20).Value = 24

# Row 127
$ws.Cells.Item(127, 4).Value = 44624
$ws.Cells.Item(127, 11).Value = "Sutil De Gase"
$ws.Cells.Item(127, 12).Value = "Primera"
$ws.Cells.Item(127, 13).Value = 200
$ws.Cells.Item(127, 14).Value = 46000
$ws.Cells.Item(127, 15).Value = 47000
$ws.Cells.Item(127, 16).Value = 46500
$ws.Cells.Item(127, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(127, 18).Value = "Perú"
$ws.Cells.Item(127, 19).Value = 1938
$ws.Cells.Item(127, 20).Value = 24

# Row 128
$ws.Cells.Item(128, 4).Value = 44624
$ws.Cells.Item(128, 11).Value = "Tahití"
$ws.Cells.Item(128, 12).Value = "Primera"
$ws.Cells.Item(128, 13).Value = 300
$ws.Cells.Item(128, 14).Value = 45000
$ws.Cells.Item(128, 15).Value = 46000
$ws.Cells.Item(128, 16).Value = 45500
$ws.Cells.Item(128, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(128, 18).Value = "Perú"
$ws.Cells.Item(128, 19).Value = 1896
$ws.Cells.Item(128, 20).Value = 24

# Row 129
$ws.Cells.Item(129, 4).Value = 44609
$ws.Cells.Item(129, 11).Value = "Sin especificar"
$ws.Cells.Item(129, 12).Value = "2a amarillo"
$ws.Cells.Item(129, 13).Value = 250
$ws.Cells.Item(129, 14).Value = 27000
$ws.Cells.Item(129, 15).Value = 28000
$ws.Cells.Item(129, 16).Value = 27500
$ws.Cells.Item(129, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(129, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(129, 19).Value = 1375
$ws.Cells.Item(129, 20).Value = 20

# Row 130
$ws.Cells.Item(130, 4).Value = 44286
$ws.Cells.Item(130, 11).Value = "Sin especificar"
$ws.Cells.Item(130, 12).Value = "2a plateado"
$ws.Cells.Item(130, 13).Value = 300
$ws.Cells.Item(130, 14).Value = 24000
$ws.Cells.Item(130, 15).Value = 25000
$ws.Cells.Item(130, 16).Value = 24500
$ws.Cells.Item(130, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(130, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(130, 19).Value = 1225
$ws.Cells.Item(130, 20).Value = 20

# Row 131
$ws.Cells.Item(131, 4).Value = 44216
$ws.Cells.Item(131, 11).Value = "Sin especificar"
$ws.Cells.Item(131, 12).Value = "1a amarillo"
$ws.Cells.Item(131, 13).Value = 250
$ws.Cells.Item(131, 14).Value = 29000
$ws.Cells.Item(131, 15).Value = 30000
$ws.Cells.Item(131, 16).Value = 29500
$ws.Cells.Item(131, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(131, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(131, 19).Value = 1475
$ws.Cells.Item(131, 20).Value = 20

# Row 132
$ws.Cells.Item(132, 4).Value = 44174
$ws.Cells.Item(132, 11).Value = "Sin especificar"
$ws.Cells.Item(132, 12).Value = "2a amarillo"
$ws.Cells.Item(132, 13).Value = 300
$ws.Cells.Item(132, 14).Value = 15000
$ws.Cells.Item(132, 15).Value = 16000
$ws.Cells.Item(132, 16).Value = 15500
$ws.Cells.Item(132, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(132, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(132, 19).Value = 775
$ws.Cells.Item(132, 20).Value = 20

# Row 133
$ws.Cells.Item(133, 4).Value = 44475
$ws.Cells.Item(133, 11).Value = "Sin especificar"
$ws.Cells.Item(133, 12).Value = "1a amarillo"
$ws.Cells.Item(133, 13).Value = 300
$ws.Cells.Item(133, 14).Value = 13000
$ws.Cells.Item(133, 15).Value = 14000
$ws.Cells.Item(133, 16).Value = 13500
$ws.Cells.Item(133, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(133, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(133, 19).Value = 675
$ws.Cells.Item(133, 20).Value = 20

# Row 134
$ws.Cells.Item(134, 4).Value = 44327
$ws.Cells.Item(134, 11).Value = "Sin especificar"
$ws.Cells.Item(134, 12).Value = "2a amarillo"
$ws.Cells.Item(134, 13).Value = 200
$ws.Cells.Item(134, 14).Value = 18000
$ws.Cells.Item(134, 15).Value = 20000
$ws.Cells.Item(134, 16).Value = 19000
$ws.Cells.Item(134, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(134, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(134, 19).Value = 950
$ws.Cells.Item(134, 20).Value = 20
